$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (columns D, M, N, O, P, S) after the reshuffle
$data = @{
    2 = @{ D = 44229; M = 55; N = 11000; O = 12000; P = 11364; S = 812 }
    3 = @{ D = 44253; M = 90; N = 12000; O = 13000; P = 12667; S = 905 }
    4 = @{ D = 44232; M = 60; N = 11000; O = 12000; P = 11583; S = 827 }
    5 = @{ D = 44172; M = 90; N = 8500;  O = 9000;  P = 8806;  S = 629 }
    6 = @{ D = 44216; M = 55; N = 11000; O = 12000; P = 11545; S = 825 }
    7 = @{ D = 44210; M = 70; N = 10000; O = 11000; P = 10357; S = 740 }
    8 = @{ D = 44181; M = 65; N = 9000;  O = 10000; P = 9462;  S = 676 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
